# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Both sheets contain identical data, and both need the same updates.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1021
    6  = 4558
    8  = 371
    9  = 1317
    10 = 543
    12 = 916
    14 = 512
    16 = 240
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
